$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updated values
$ws.Range("B2").Value = [double]"8.427485376216737e-09"
$ws.Range("C2").Value = [double]"2.220651329265522e-06"
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("G2").Value = 7.194163157364425

# Row 3 updated values
$ws.Range("B3").Value = 0.02258322285507441
$ws.Range("C3").Value = 0.05231270169004087
$ws.Range("D3").Value = 0.1529057820181812
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 0.7276884136373534
